$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: grow the row height to fit the new wrapped text ---
$ws.Rows.Item(24).RowHeight = 91

# --- Row 25: new summary row (90% threshold) ---
$ws.Rows.Item(25).RowHeight = 58
$ws.Range("B25:F25").Merge()
$ws.Range("B25").Value = " Min. 1st Qu.  Median    Mean 3rd Qu.    Max. `n     12      22      29      32      39      91"
$ws.Range("G24").Value = "79/100"
$ws.Range("G25").Value = "90/91"

# --- Row 26: new summary row (80% threshold) ---
$ws.Rows.Item(26).RowHeight = 43
$ws.Range("B26:F26").Merge()
$ws.Range("B26").Value = " Min. 1st Qu.  Median    Mean 3rd Qu.    Max. `n   9.00   15.00   20.00   22.17   26.00   55.00"
$ws.Range("G26").Value = "97/100"

# --- Row 27: new summary row (75% threshold) ---
$ws.Rows.Item(27).RowHeight = 44
$ws.Range("B27:F27").Merge()
$ws.Range("G27").Value = "96/100"
$ws.Range("B27").Value = " Min. 1st Qu.  Median    Mean 3rd Qu.    Max. `n   7.00   14.00   17.00   19.28   22.00   61.00 "

# --- Row 28: new summary row (70% threshold) ---
$ws.Rows.Item(28).RowHeight = 49
$ws.Range("B28:F28").Merge()
$ws.Range("G28").Value = "90/100"
$ws.Range("B28").Value = " Min. 1st Qu.  Median    Mean 3rd Qu.    Max. `n   8.00   12.00   15.00   16.54   20.00   36.00 "

# --- numeric values + number format for column A (percent thresholds) ---
$ws.Range("A25").Value = 0.9
$ws.Range("A26").Value = 0.8
$ws.Range("A27").Value = 0.75
$ws.Range("A28").Value = 0.7
$ws.Range("A25:A28").NumberFormat = "0%"

# --- alignment / wrap formatting to mirror the rows above (18-24) ---
$ws.Range("B25:F25").HorizontalAlignment = -4108
$ws.Range("B25:F25").WrapText = $true

$ws.Range("B26").HorizontalAlignment = -4108
$ws.Range("B26").WrapText = $true
$ws.Range("C26:F26").HorizontalAlignment = -4108

$ws.Range("B27").HorizontalAlignment = -4108
$ws.Range("B27").WrapText = $true
$ws.Range("C27:F27").HorizontalAlignment = -4108

$ws.Range("B28").HorizontalAlignment = -4108
$ws.Range("B28").WrapText = $true
$ws.Range("C28:F28").HorizontalAlignment = -4108

# --- view / selection changes ---
$ws.Range("B28:F28").Select()

# --- page setup ---
$ws.PageSetup.Orientation = 1

Write-Host "edit complete"
